# Updates the "cryptos" price/volume sheet with freshly scraped values
# (mirrors the automated "Updated cryptos list ... with GitHub Actions" commit).
#
# Every Price/Volume(1h) cell is stored as plain text in the workbook (inlineStr),
# even though many values look numeric (e.g. "0.429", "10.60", "0.0236"). To keep
# them as text (and not let Excel silently coerce them to numbers / drop
# significant trailing zeros), cells whose new value would otherwise be
# auto-recognized as a number are temporarily forced to a Text number format
# before the value is written, then restored to the default "Normal" style
# afterwards so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.542.24'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '2.598.08'
$ws.Range('E3').Value = '  +1.73%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '515.39'
$ws.Range('E5').Value = '  +2.62%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.53'
$ws.Range('E6').Value = '  +1.73%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.601'
$ws.Range('E8').Value = '  +4.68%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.64'
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('E10').Value = '  +2.38%  '
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('D13').Value = '3.055.49'
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('D14').Value = '60.590.69'
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.81'
$ws.Range('E15').Value = '  +1.97%  '
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').Value = '2.608.06'
$ws.Range('E17').Value = '  +1.45%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '4.76'
$ws.Range('E18').Value = '  -0.35%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '353.23'
$ws.Range('E19').Value = '  +2.44%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.60'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.23'
$ws.Range('E21').Value = '  +3.69%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '61.05'
$ws.Range('E23').Value = '  +2.13%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.429'
$ws.Range('E24').Value = '  +3.20%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').Value = '2.720.21'
$ws.Range('E25').Value = '  +1.53%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.166'
$ws.Range('E26').Value = '  +1.02%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').Value = '0.0₃0841'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.35'
$ws.Range('E29').Value = '  -0.73%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.21'
$ws.Range('E31').Value = '  +9.40%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '19.44'
$ws.Range('E32').Value = '  +1.91%  '
$ws.Range('E33').Value = '  +2.83%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '149.74'
$ws.Range('E34').Value = '  -3.59%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.07'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.943'
$ws.Range('E36').Value = '  +10.66%  '
$ws.Range('E37').Value = '  +1.33%  '
$ws.Range('E38').Value = '  +2.38%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '36.33'
$ws.Range('E39').Value = '  +2.61%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.77'
$ws.Range('E40').Value = '  +0.78%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.841'
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '286.65'
$ws.Range('E42').Value = '  -2.97%  '
$ws.Range('E43').Value = '  +2.68%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.626'
$ws.Range('E44').Value = '  +2.39%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0559'
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.96'
$ws.Range('E47').Value = '  +2.60%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '19.52'
$ws.Range('E48').Value = '  +0.15%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0236'
$ws.Range('E49').Value = '  +1.72%  '
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').Value = '1.977.56'
$ws.Range('E51').Value = '  -0.69%  '

# Restore default (Normal) style on cells where we forced a Text number format,
# so no extra explicit style index lingers on these cells.
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
